$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Опросник")
$ws.Columns.Item(1).Delete()
